$d = $word.ActiveDocument

# ============================================================================
# Change 1: "(3) correlated" -> "(3) " + strikethrough "correlated" +
#           bold-red " endoswitch"
# ============================================================================
$rng = $d.Content
$rng.Find.Execute("correlated", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$correlatedStart = $rng.Start
$correlatedEnd = $rng.End

# strike the existing word "correlated"
$strikeRng = $d.Range($correlatedStart, $correlatedEnd)
$strikeRng.Font.StrikeThrough = $true

# insert " endoswitch" right after it
$insertPoint = $d.Range($correlatedEnd, $correlatedEnd)
$insertPoint.InsertAfter(" endoswitch")

# format the leading space and the new word bold+red (applied separately so
# the space/word keep their own run, matching the target run layout)
$spaceRng = $d.Range($correlatedEnd, $correlatedEnd + 1)
$spaceRng.Font.Bold = $true
$spaceRng.Font.Color = 255

$wordRng = $d.Range($correlatedEnd + 1, $correlatedEnd + 11)
$wordRng.Font.Bold = $true
$wordRng.Font.Color = 255

# ============================================================================
# Change 2: "...model with exogenous switching" (the occurrence right before
#           "Number of observations =") -> "...model with endogenous switching"
# ============================================================================
$full = $d.Content.Text
$anchor = $full.IndexOf("Number of observations")
$search = "exogenous switching"
$pos = $full.LastIndexOf($search, $anchor)

$xStart = $pos + 1          # the lone "x" run
$xEnd = $pos + 2
$xRng = $d.Range($xStart, $xEnd)
$xRng.Text = "ndog"

# the run originally reading "ogenous switching" now starts 3 chars later
# (it grew by 3: "x" -> "ndog"); trim its leading "og" so the combined text
# reads "...with endogenous switching"
$ogStart = $xEnd + 3
$ogRng = $d.Range($ogStart, $ogStart + 2)
$ogRng.Text = ""

# ============================================================================
# Change 3 & 4: move the "_GoBack" bookmark from the "y5_positive ... (+)"
#               line down to the "rho p ... rho(+)" results line (right
#               before "rho" through just after the trailing space of
#               "rho(+) ").
# ============================================================================
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks("_GoBack").Delete()
}

$full2 = $d.Content.Text
$rhoParaPos = $full2.IndexOf("rho p")
$rhoPlusPos = $full2.IndexOf("rho(+)", $rhoParaPos)
$afterRhoPlus = $rhoPlusPos + 6   # length of "rho(+)"
$bmEnd = $afterRhoPlus + 1        # just after the trailing space

$bmRange = $d.Range($rhoParaPos, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
$goBack = $d.Bookmarks("_GoBack")
$goBack.Start = $bmEnd
$goBack.End = $bmEnd

"done"
